# Auto-generated edit script: apply artfynd row-data corrections for rows 10-22, 24-28
# (field values were re-matched/corrected between observation records; see commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 112044191
$ws.Range("B10").Value = 96735
$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value = 554719
$ws.Range("R10").Value = 6698669
$ws.Range("AF10").ClearContents()
$ws.Range("A11").Value = 112044188
$ws.Range("B11").Value = 96735
$ws.Range("Q11").Value = 554647
$ws.Range("R11").Value = 6698760
$ws.Range("A12").Value = 112044189
$ws.Range("B12").Value = 96735
$ws.Range("Q12").Value = 554686
$ws.Range("R12").Value = 6698721
$ws.Range("A13").Value = 112044193
$ws.Range("B13").Value = 96735
$ws.Range("Q13").Value = 554737
$ws.Range("R13").Value = 6698616
$ws.Range("A14").Value = 112044185
$ws.Range("B14").Value = 96735
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = "Knärot"
$ws.Range("G14").Value = "Goodyera repens"
$ws.Range("H14").Value = "(L.) R. Br."
$ws.Range("Q14").Value = 554752
$ws.Range("R14").Value = 6698637
$ws.Range("A15").Value = 112044194
$ws.Range("B15").Value = 96735
$ws.Range("Q15").Value = 554746
$ws.Range("R15").Value = 6698619
$ws.Range("A16").Value = 112044187
$ws.Range("B16").Value = 96735
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("Q16").Value = 554629
$ws.Range("R16").Value = 6698775
$ws.Range("A17").Value = 112044192
$ws.Range("B17").Value = 96735
$ws.Range("Q17").Value = 554727
$ws.Range("R17").Value = 6698622
$ws.Range("A18").Value = 112044170
$ws.Range("B18").Value = 89993
$ws.Range("E18").Value = 1209
$ws.Range("F18").Value = "Rynkskinn"
$ws.Range("G18").Value = "Phlebia centrifuga"
$ws.Range("H18").Value = "P.Karst."
$ws.Range("Q18").Value = 554745
$ws.Range("R18").Value = 6698641
$ws.Range("A19").Value = 112044190
$ws.Range("B19").Value = 96735
$ws.Range("Q19").Value = 554682
$ws.Range("R19").Value = 6698694
$ws.Range("A20").Value = 112044186
$ws.Range("B20").Value = 96735
$ws.Range("Q20").Value = 554675
$ws.Range("R20").Value = 6698785
$ws.Range("A21").Value = 112044174
$ws.Range("B21").Value = 96735
$ws.Range("Q21").Value = 554690
$ws.Range("R21").Value = 6698722
$ws.Range("A22").Value = 112044163
$ws.Range("B22").Value = 56575
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 103021
$ws.Range("F22").Value = "Talltita"
$ws.Range("G22").Value = "Poecile montanus"
$ws.Range("H22").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q22").Value = 554650
$ws.Range("R22").Value = 6698762
$ws.Range("A24").Value = 112044184
$ws.Range("B24").Value = 96735
$ws.Range("D24").Value = "VU"
$ws.Range("E24").Value = 220787
$ws.Range("F24").Value = "Knärot"
$ws.Range("G24").Value = "Goodyera repens"
$ws.Range("H24").Value = "(L.) R. Br."
$ws.Range("Q24").Value = 554833
$ws.Range("R24").Value = 6698646
$ws.Range("A25").Value = 112044195
$ws.Range("B25").Value = 96735
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("Q25").Value = 554806
$ws.Range("R25").Value = 6698598
$ws.Range("A26").Value = 112044158
$ws.Range("B26").Value = 89553
$ws.Range("E26").Value = 1202
$ws.Range("F26").Value = "Ullticka"
$ws.Range("G26").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H26").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K26").ClearContents()
$ws.Range("L26").ClearContents()
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("Q26").Value = 554756
$ws.Range("R26").Value = 6698631
$ws.Range("A27").Value = 112044171
$ws.Range("B27").Value = 89834
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 658
$ws.Range("F27").Value = "Rosenticka"
$ws.Range("G27").Value = "Rhodofomes roseus"
$ws.Range("H27").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q27").Value = 554758
$ws.Range("R27").Value = 6698625
$ws.Range("A28").Value = 112044162
$ws.Range("B28").Value = 56575
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 103021
$ws.Range("F28").Value = "Talltita"
$ws.Range("G28").Value = "Poecile montanus"
$ws.Range("H28").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q28").Value = 554765
$ws.Range("R28").Value = 6698666
$ws.Range("J21").Style = "Normal"
$ws.Range("K21").Value = "överblommad"
$ws.Range("L21").Style = "Normal"
$ws.Range("N21").Style = "Normal"
$ws.Range("AF21").Style = "Normal"
$ws.Range("K22").Style = "Normal"
$ws.Range("L22").Style = "Normal"
$ws.Range("M22").Value = "spel/sång"
$ws.Range("N22").Style = "Normal"
$ws.Range("K28").Style = "Normal"
$ws.Range("L28").Style = "Normal"
$ws.Range("M28").Value = "spel/sång"
$ws.Range("N28").Style = "Normal"
